$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N to make room for the
# "Variable Instalments" RBI column — this shifts the existing
# N:O (Late/Outstanding) columns right to O:Q.
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with R6 selected,
# matching the new workbook view state.
$ws.Activate() | Out-Null
$ws.Range("R6").Select() | Out-Null
